# Update data files - Bot run at 2026-02-20 21:30:48 UTC
# Apply the numeric updates to row 2 (llama-3.1-8b-instant) of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 0.0004861111111111111
$ws.Range("K2").Value = 3331
$ws.Range("L2").Value = 0.006662
